# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# OFF sheet - Home row (row 2) cumulative target depth stats
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 315
$wsOff.Range("C2").Value = 227
$wsOff.Range("D2").Value = 169
$wsOff.Range("E2").Value = 78
$wsOff.Range("G2").Value = 9

# DEF sheet - Home row (row 2) cumulative target depth stats
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 364
$wsDef.Range("C2").Value = 251
$wsDef.Range("D2").Value = 86
$wsDef.Range("E2").Value = 41
